# Update crypto price (column D) and volume change percentage (column E)
# values on the active sheet, row by row, to match the latest scrape.
#
# Prices in column D are stored as plain text (e.g. thousand separators use
# '.' rather than ',', which Excel would otherwise try to parse as a
# number/date). We prefix values with a leading apostrophe so Excel COM
# keeps them as literal text instead of coercing them to numbers, matching
# the original inline-string cell content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @(newPriceOrNull, newVolumePctOrNull)
# A $null entry means that column is left untouched for that row.
$updates = @{
    2  = @('23.793.36', '  +1.88%  ')
    3  = @('1.656.59', '  +1.80%  ')
    4  = @($null, '  -0.01%  ')
    5  = @($null, '  +0.08%  ')
    6  = @('304.16', '  +0.75%  ')
    7  = @('0.3807', $null)
    8  = @('0.3634', '  +0.23%  ')
    9  = @('50.95', '  -1.13%  ')
    10 = @('1.259', '  +3.46%  ')
    11 = @('0.08237', '  +0.85%  ')
    12 = @($null, '  -0.01%  ')
    13 = @('22.75', '  +2.60%  ')
    14 = @('6.554', '  +1.51%  ')
    15 = @('7.488', '  +2.76%  ')
    16 = @('0.00001243', '  +0.50%  ')
    17 = @('1.655.75', '  +3.41%  ')
    18 = @('97.86', '  +3.42%  ')
    19 = @('0.06993', '  +0.66%  ')
    20 = @('6.817', '  +4.04%  ')
    21 = @('17.82', '  +1.68%  ')
    22 = @($null, '  +0.05%  ')
    23 = @('12.88', '  +3.08%  ')
    24 = @('23.781.39', '  +1.81%  ')
    25 = @('2.555', $null)
    26 = @('3.087', '  +0.64%  ')
    27 = @('21.37', '  +1.17%  ')
    28 = @('151.69', '  +0.87%  ')
    29 = @('5.236', '  -0.54%  ')
    30 = @('134.41', '  +1.17%  ')
    31 = @('1.839.68', '  +2.64%  ')
    32 = @('6.931', '  +5.07%  ')
    33 = @($null, '  +1.75%  ')
    34 = @('1.081', '  +1.96%  ')
    35 = @('11.88', '  +6.25%  ')
    36 = @('0.02840', '  +3.05%  ')
    37 = @('0.2534', '  +1.85%  ')
    38 = @('6.152', '  +3.26%  ')
    39 = @('0.08831', '  +0.65%  ')
    40 = @('0.07128', '  +0.08%  ')
    41 = @('13.29', '  +11.11%  ')
    42 = @('0.7093', '  +1.82%  ')
    43 = @('1.347', '  +1.57%  ')
    44 = @('16.09', '  +1.74%  ')
    45 = @('0.6576', '  +2.22%  ')
    46 = @('2.340', '  +3.17%  ')
    47 = @($null, '  +0.07%  ')
    48 = @('3.967', '  +0.30%  ')
    49 = @('0.07966', $null)
    50 = @('128.43', '  +1.14%  ')
    51 = @('1.199', '  +1.03%  ')
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $newPrice = $vals[0]
    $newVolume = $vals[1]

    if ($null -ne $newPrice) {
        $ws.Cells.Item($row, 4).Value = "'" + $newPrice
    }
    if ($null -ne $newVolume) {
        $ws.Cells.Item($row, 5).Value = $newVolume
    }
}
